$wb = $excel.ActiveWorkbook

# --- Productdata sheet: divide InventoryCosts (D), BackorderCosts (F) and LostSale (I)
#     by 2500 for rows 2-11 ---
$ws = $wb.Worksheets.Item("Productdata")

for ($r = 2; $r -le 11; $r++) {
    $cellD = $ws.Cells.Item($r, 4)
    $cellD.Value = $cellD.Value() / 2500   # column D

    $cellF = $ws.Cells.Item($r, 6)
    $cellF.Value = $cellF.Value() / 2500   # column F

    $cellI = $ws.Cells.Item($r, 9)
    $cellI.Value = $cellI.Value() / 2500   # column I
}

# --- ForcastedStandardDeviation sheet: zero out columns B-E for rows 9-11 ---
$ws2 = $wb.Worksheets.Item("ForcastedStandardDeviation")

$ws2.Range("B9:E11").Value = 0
